$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 410.25757
$ws.Range("I17").Value = 700
$ws.Range("J17").Value = 386.5082
$ws.Range("K17").Value = 2100
$ws.Range("L17").Value = 1159.5246
$ws.Range("M17").Value = -1932
$ws.Range("N17").Value = -1495.5246
$ws.Range("H33").Value = 6191.6665
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 6191.6665
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 6191.6665
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -6649.6665
$ws.Range("H112").Value = 3810.9546
$ws.Range("J112").Value = 3892.6047
$ws.Range("L112").Value = 11677.8141
$ws.Range("N112").Value = -13893.8141
$ws.Range("H129").Value = 1025.6296
$ws.Range("J129").Value = 1181.591
$ws.Range("L129").Value = 3544.773
$ws.Range("N129").Value = -13544.773
$ws.Range("H132").Value = 36523.32
$ws.Range("I132").Value = 39284.117
$ws.Range("J132").Value = 633
$ws.Range("K132").Value = 117852.351
$ws.Range("L132").Value = 1899
$ws.Range("M132").Value = -115322.351
$ws.Range("N132").Value = -6959
$ws.Range("H137").Value = 951.2222
$ws.Range("I137").Value = 832.381
$ws.Range("K137").Value = 2497.143
$ws.Range("M137").Value = 52.85699999999997
$ws.Range("H141").Value = 3648.7693
$ws.Range("J141").Value = 3580
$ws.Range("L141").Value = 10740
$ws.Range("N141").Value = -21100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 359569.22
$ws.Range("I32").Value = 3766.9265
$ws.Range("J32").Value = 2220689
$ws.Range("K32").Value = 3766.9265
$ws.Range("L32").Value = 2220689
$ws.Range("M32").Value = -3479.9265
$ws.Range("N32").Value = -2221263
$ws.Range("H61").Value = 1377.8667
$ws.Range("I61").Value = 1232.2222
$ws.Range("J61").Value = 1596.3334
$ws.Range("K61").Value = 1232.2222
$ws.Range("L61").Value = 1596.3334
$ws.Range("M61").Value = -1020.2222
$ws.Range("N61").Value = -2020.3334
$ws.Range("H74").Value = 1108.5483
$ws.Range("I74").Value = 1027.3214
$ws.Range("J74").Value = 1866.6666
$ws.Range("K74").Value = 1027.3214
$ws.Range("L74").Value = 1866.6666
$ws.Range("M74").Value = -153.3214
$ws.Range("N74").Value = -3614.6666
$ws.Range("H77").Value = 1108.5483
$ws.Range("I77").Value = 1027.3214
$ws.Range("J77").Value = 1866.6666
$ws.Range("K77").Value = 5136.607
$ws.Range("L77").Value = 9333.333000000001
$ws.Range("M77").Value = -768.607
$ws.Range("N77").Value = -18069.333
$ws.Range("H132").Value = 16967758
$ws.Range("I132").Value = 27028510
$ws.Range("K132").Value = 81085530
$ws.Range("M132").Value = -81083000
$ws.Range("H136").Value = 1377.8667
$ws.Range("I136").Value = 1232.2222
$ws.Range("J136").Value = 1596.3334
$ws.Range("K136").Value = 3696.6666
$ws.Range("L136").Value = 4789.0002
$ws.Range("M136").Value = -1146.6666
$ws.Range("N136").Value = -9889.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 20627.166
$ws.Range("I96").Value = 18709.166
$ws.Range("J96").Value = 22545.166
$ws.Range("K96").Value = 18709.166
$ws.Range("L96").Value = 22545.166
$ws.Range("M96").Value = -15963.166
$ws.Range("N96").Value = -28037.166
$ws.Range("H126").Value = 41814.547
$ws.Range("J126").Value = 41814.547
$ws.Range("L126").Value = 41814.547
$ws.Range("N126").Value = -51694.547
$ws.Range("H134").Value = 4462.396
$ws.Range("I134").Value = 1651.8788
$ws.Range("J134").Value = 10645.533
$ws.Range("K134").Value = 4955.636399999999
$ws.Range("L134").Value = 31936.599
$ws.Range("M134").Value = -2420.636399999999
$ws.Range("N134").Value = -37006.599

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1284.1875
$ws.Range("I16").Value = 1205.5
$ws.Range("J16").Value = 1415.3334
$ws.Range("K16").Value = 1205.5
$ws.Range("L16").Value = 1415.3334
$ws.Range("M16").Value = -918.5
$ws.Range("N16").Value = -1989.3334
$ws.Range("H31").Value = 12506285
$ws.Range("I31").Value = 20836374
$ws.Range("J31").Value = 11151.75
$ws.Range("K31").Value = 20836374
$ws.Range("L31").Value = 11151.75
$ws.Range("M31").Value = -20836079
$ws.Range("N31").Value = -11741.75
$ws.Range("H34").Value = 12506285
$ws.Range("I34").Value = 20836374
$ws.Range("J34").Value = 11151.75
$ws.Range("K34").Value = 20836374
$ws.Range("L34").Value = 11151.75
$ws.Range("M34").Value = -20836172
$ws.Range("N34").Value = -11555.75
$ws.Range("H99").Value = 271314.72
$ws.Range("I99").Value = 406149.47
$ws.Range("J99").Value = 1645.238
$ws.Range("K99").Value = 406149.47
$ws.Range("L99").Value = 1645.238
$ws.Range("M99").Value = -404651.47
$ws.Range("N99").Value = -4641.238
$ws.Range("H105").Value = 1346.7142
$ws.Range("I105").Value = 995.4
$ws.Range("K105").Value = 995.4
$ws.Range("M105").Value = 751.6
$ws.Range("H113").Value = 1284.1875
$ws.Range("I113").Value = 1205.5
$ws.Range("J113").Value = 1415.3334
$ws.Range("K113").Value = 1205.5
$ws.Range("L113").Value = 1415.3334
$ws.Range("M113").Value = 964.5
$ws.Range("N113").Value = -5755.3334
$ws.Range("H126").Value = 271314.72
$ws.Range("I126").Value = 406149.47
$ws.Range("J126").Value = 1645.238
$ws.Range("K126").Value = 1218448.41
$ws.Range("L126").Value = 4935.714
$ws.Range("M126").Value = -1215978.41
$ws.Range("N126").Value = -9875.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3218.1333
$ws.Range("I7").Value = 1999.6666
$ws.Range("K7").Value = 1999.6666
$ws.Range("M7").Value = -1887.6666
$ws.Range("H40").Value = 2238.5862
$ws.Range("I40").Value = 2071.3333
$ws.Range("K40").Value = 2071.3333
$ws.Range("M40").Value = -1935.3333
$ws.Range("H92").Value = 21880
$ws.Range("J92").Value = 21880
$ws.Range("L92").Value = 21880
$ws.Range("N92").Value = -26872
$ws.Range("H126").Value = 3218.1333
$ws.Range("I126").Value = 1999.6666
$ws.Range("K126").Value = 5998.9998
$ws.Range("M126").Value = -3528.9998
$ws.Range("H136").Value = 8932.044
$ws.Range("I136").Value = 8949.941000000001
$ws.Range("J136").Value = 8881.333000000001
$ws.Range("K136").Value = 26849.823
$ws.Range("L136").Value = 26643.999
$ws.Range("M136").Value = -24299.823
$ws.Range("N136").Value = -31743.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 62501716
$ws.Range("J126").Value = 2028.75
$ws.Range("L126").Value = 6086.25
$ws.Range("N126").Value = -11026.25
$ws.Range("H136").Value = 67927.60000000001
$ws.Range("I136").Value = 78108.766
$ws.Range("K136").Value = 234326.298
$ws.Range("M136").Value = -231776.298
